$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("measurements")

# Give the new cells the same formatting as the rest of the data row (A5)
# before filling them in, so they share its existing cell style.
$ws.Range("A5").Copy()
$ws.Range("B5:C5").PasteSpecial(-4122)

# Populate the new data cells next to the existing A5 value.
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "A"

# Remove the extra empty row 6 (A6), which is no longer needed.
$ws.Rows.Item(6).Delete()

# Register the newly-used columns B and C at the sheet's standard width.
$ws.Columns.Item(2).ColumnWidth = $ws.StandardWidth
$ws.Columns.Item(3).ColumnWidth = $ws.StandardWidth
